$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "El cuarto de Jacob"
$ws.Range("B5").Value = "Virginia Woolf"
$ws.Range("C5").Value = "Lumen"
